# Updated capital structure database
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "'3"
$ws.Range("K2").Value = -1.734
$ws.Range("U2").Value = 0.17
$ws.Range("V2").Value = 0.005407124681933843
$ws.Range("W2").Value = -0.04872549019607843
$ws.Range("X2").Value = 0.05981543109285033
$ws.Range("Y2").Value = -0.1085409212889288
$ws.Range("AA2").Value = -0.04757760371337395
$ws.Range("AB2").Value = 0.05817469445585257
$ws.Range("AC2").Value = -0.1058854123940585
$ws.Range("AD2").Value = 1.47
$ws.Range("AF2").Value = 1.47
$ws.Range("AG2").Value = 1.3
$ws.Range("AH2").Value = 0.04466727438468551
$ws.Range("AI2").Value = 0.03611793611793612
$ws.Range("AJ2").Value = 0.03970678069639585
$ws.Range("AK2").Value = 0.03207500616827041
$ws.Range("AL2").Value = 0.003
$ws.Range("AM2").Value = 0.003
$ws.Range("AN2").Value = -1.233221476510067
$ws.Range("AO2").Value = -576.3333333333334
$ws.Range("AP2").Value = -1.090604026845638
$ws.Range("AQ2").Value = -576.3333333333334

# Row 3
$ws.Range("K3").Value = -0.753
$ws.Range("U3").Value = 0.13
$ws.Range("V3").Value = 0.007027027027027027
$ws.Range("W3").Value = -0.0378391959798995
$ws.Range("X3").Value = 0.06005851874003106
$ws.Range("Y3").Value = -0.09789771471993056
$ws.Range("AA3").Value = -0.03674963396778917
$ws.Range("AB3").Value = 0.05817469445585257
$ws.Range("AC3").Value = -0.09492432842364174
$ws.Range("AD3").Value = 1.2
$ws.Range("AF3").Value = 1.2
$ws.Range("AG3").Value = 1.07
$ws.Range("AH3").Value = 0.06091370558375635
$ws.Range("AI3").Value = 0.05769230769230769
$ws.Range("AJ3").Value = 0.05467552376085845
$ws.Range("AK3").Value = 0.0517658442186744
$ws.Range("AN3").Value = -1.6
$ws.Range("AP3").Value = -1.426666666666667

# Row 4
$ws.Range("B4").Value = 'Karelian Diamond Resources Plc (AIM:KDR)'
$ws.Range("K4").Value = -0.497
$ws.Range("U4").Value = 0.018
$ws.Range("V4").Value = 0.003938730853391684
$ws.Range("W4").Value = -0.04872549019607843
$ws.Range("X4").Value = 0.05981543109285033
$ws.Range("Y4").Value = -0.1085409212889288
$ws.Range("AA4").Value = -0.04757760371337395
$ws.Range("AB4").Value = 0.05830780868068455
$ws.Range("AC4").Value = -0.1058854123940585
$ws.Range("AD4").Value = 0.27
$ws.Range("AF4").Value = 0.27
$ws.Range("AG4").Value = 0.252
$ws.Range("AH4").Value = 0.05578512396694216
$ws.Range("AI4").Value = 0.02603664416586307
$ws.Range("AJ4").Value = 0.05226047283284944
$ws.Range("AK4").Value = 0.02434312210200927
$ws.Range("AL4").Value = 0.003
$ws.Range("AM4").Value = 0.003
$ws.Range("AO4").Value = -164
$ws.Range("AQ4").Value = -164

# Row 5
$ws.Range("K5").Value = -0.484
$ws.Range("U5").Value = 0.022
$ws.Range("V5").Value = 0.002628434886499403
$ws.Range("W5").Value = -0.04953940634595701
$ws.Range("X5").Value = 0.0573323578869859
$ws.Range("Y5").Value = -0.1068717642329429
$ws.Range("AA5").Value = -0.04963084495488106
$ws.Range("AB5").Value = 0.0573323578869859
$ws.Range("AC5").Value = -0.106963202841867
$ws.Range("AG5").Value = -0.022
$ws.Range("AJ5").Value = -0.002635361763296598
$ws.Range("AK5").Value = -0.002313840976020194
$ws.Range("AP5").Value = 0.04977375565610859

# Row 6 (Kibo Energy PLC) was removed from the dataset entirely
$ws.Rows(6).Delete()

